$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.926332950592041
$ws.Range("B1").Value = 3.47041392326355
$ws.Range("C1").Value = 2.62143087387085
$ws.Range("D1").Value = 2.041930437088013
$ws.Range("E1").Value = 2.314959526062012
